$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dades_Meteo")

$ws.Range("E2").Value = "2026-02-13 19:48:48"
$ws.Range("M2").Value = "1.3 °C 19:18 TU"
$ws.Range("O2").Value = "-0.6 °C"
$ws.Range("E3").Value = "2026-02-13 19:48:51"
$ws.Range("I3").Value = "5.9 mm"
$ws.Range("E4").Value = "2026-02-13 19:48:54"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "76%"
$ws.Range("H4").NumberFormat = "general"
$ws.Range("I4").Value = "6.2 mm"
$ws.Range("J4").Value = "994.5 hPa"
$ws.Range("E5").Value = "2026-02-13 19:48:56"
$ws.Range("G5").Value = "106 cm"
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = "82%"
$ws.Range("H5").NumberFormat = "general"
$ws.Range("I5").Value = "1.0 mm"
$ws.Range("E6").Value = "2026-02-13 19:48:59"
$ws.Range("I6").Value = "5.1 mm"
$ws.Range("J6").Value = "994.5 hPa"
$ws.Range("O6").Value = "9.1 °C"
$ws.Range("E7").Value = "2026-02-13 19:49:02"
$ws.Range("J7").Value = "994.8 hPa"
$ws.Range("L7").Value = "55.1 km/h - 277º 19:20 TU"
$ws.Range("E8").Value = "2026-02-13 19:49:04"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "80%"
$ws.Range("H8").NumberFormat = "general"
$ws.Range("J8").Value = "994.7 hPa"
$ws.Range("O8").Value = "9.3 °C"
$ws.Range("E9").Value = "2026-02-13 19:49:07"
$ws.Range("I9").Value = "4.0 mm"
$ws.Range("E10").Value = "2026-02-13 19:49:10"
$ws.Range("I10").Value = "19.1 mm"
$ws.Range("E11").Value = "2026-02-13 19:49:12"
$ws.Range("E12").Value = "2026-02-13 19:49:15"
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = "85%"
$ws.Range("H12").NumberFormat = "general"
$ws.Range("I12").Value = "5.5 mm"
$ws.Range("E13").Value = "2026-02-13 19:49:17"
$ws.Range("E14").Value = "2026-02-13 19:49:20"
$ws.Range("E15").Value = "2026-02-13 19:49:23"
$ws.Range("I15").Value = "4.5 mm"
$ws.Range("O15").Value = "9.6 °C"
$ws.Range("E16").Value = "2026-02-13 19:49:25"
$ws.Range("E17").Value = "2026-02-13 19:49:28"
$ws.Range("H17").NumberFormat = "@"
$ws.Range("H17").Value = "90%"
$ws.Range("H17").NumberFormat = "general"
$ws.Range("I17").Value = "5.4 mm"
$ws.Range("N17").Value = "-0.9 °C 18:53 TU"
$ws.Range("E18").Value = "2026-02-13 19:49:31"
$ws.Range("H18").NumberFormat = "@"
$ws.Range("H18").Value = "83%"
$ws.Range("H18").NumberFormat = "general"
$ws.Range("I18").Value = "10.4 mm"
$ws.Range("J18").Value = "994.7 hPa"
$ws.Range("E19").Value = "2026-02-13 19:49:33"
$ws.Range("O19").Value = "3.8 °C"
$ws.Range("E20").Value = "2026-02-13 19:49:36"
$ws.Range("I20").Value = "22.5 mm"
$ws.Range("E21").Value = "2026-02-13 19:49:39"
$ws.Range("J21").Value = "997.6 hPa"
$ws.Range("O21").Value = "1.0 °C"
$ws.Range("E22").Value = "2026-02-13 19:49:41"
$ws.Range("G22").Value = "136 cm"
$ws.Range("L22").Value = "41.8 km/h - 334º 19:26 TU"
$ws.Range("E23").Value = "2026-02-13 19:49:44"
$ws.Range("G23").Value = "188 cm"
$ws.Range("H23").NumberFormat = "@"
$ws.Range("H23").Value = "83%"
$ws.Range("H23").NumberFormat = "general"
$ws.Range("I23").Value = "9.7 mm"
$ws.Range("E24").Value = "2026-02-13 19:49:47"
$ws.Range("J24").Value = "995.5 hPa"
$ws.Range("E25").Value = "2026-02-13 19:49:50"
$ws.Range("I25").Value = "9.0 mm"
$ws.Range("L25").Value = "49.3 km/h - 260º 19:04 TU"
$ws.Range("E26").Value = "2026-02-13 19:49:53"
$ws.Range("E27").Value = "2026-02-13 19:49:55"
$ws.Range("E28").Value = "2026-02-13 19:49:58"
$ws.Range("J28").Value = "994.9 hPa"
$ws.Range("E29").Value = "2026-02-13 19:50:01"
$ws.Range("I29").Value = "13.7 mm"
$ws.Range("E30").Value = "2026-02-13 19:50:03"
$ws.Range("I30").Value = "5.6 mm"
$ws.Range("J30").Value = "994.4 hPa"
$ws.Range("E31").Value = "2026-02-13 19:50:06"
$ws.Range("H31").NumberFormat = "@"
$ws.Range("H31").Value = "73%"
$ws.Range("H31").NumberFormat = "general"
$ws.Range("I31").Value = "4.2 mm"
$ws.Range("J31").Value = "993.4 hPa"
$ws.Range("L31").Value = "100.8 km/h - 8º 19:11 TU"
$ws.Range("O31").Value = "10.3 °C"
$ws.Range("E32").Value = "2026-02-13 19:50:09"
$ws.Range("L32").Value = "44.6 km/h - 310º 19:14 TU"
$ws.Range("E33").Value = "2026-02-13 19:50:12"
$ws.Range("J33").Value = "996.4 hPa"
$ws.Range("O33").Value = "1.1 °C"
$ws.Range("E34").Value = "2026-02-13 19:50:14"
$ws.Range("G34").Value = "110 cm"
$ws.Range("E35").Value = "2026-02-13 19:50:17"
$ws.Range("H35").NumberFormat = "@"
$ws.Range("H35").Value = "76%"
$ws.Range("H35").NumberFormat = "general"
$ws.Range("J35").Value = "995.5 hPa"
$ws.Range("O35").Value = "6.0 °C"
$ws.Range("E36").Value = "2026-02-13 19:50:20"
$ws.Range("I36").Value = "8.8 mm"
$ws.Range("J36").Value = "994.5 hPa"
$ws.Range("L36").Value = "34.6 km/h - 319º 19:28 TU"
$ws.Range("E37").Value = "2026-02-13 19:50:23"
$ws.Range("J37").Value = "996.4 hPa"
$ws.Range("E38").Value = "2026-02-13 19:50:25"
$ws.Range("I38").Value = "14.1 mm"
$ws.Range("E39").Value = "2026-02-13 19:50:28"
$ws.Range("H39").NumberFormat = "@"
$ws.Range("H39").Value = "79%"
$ws.Range("H39").NumberFormat = "general"
$ws.Range("I39").Value = "19.3 mm"
$ws.Range("E40").Value = "2026-02-13 19:50:31"
$ws.Range("J40").Value = "998.0 hPa"
$ws.Range("E41").Value = "2026-02-13 19:50:33"
$ws.Range("J41").Value = "994.9 hPa"
$ws.Range("E42").Value = "2026-02-13 19:50:36"
$ws.Range("I42").Value = "10.0 mm"
$ws.Range("E43").Value = "2026-02-13 19:50:39"
$ws.Range("H43").NumberFormat = "@"
$ws.Range("H43").Value = "88%"
$ws.Range("H43").NumberFormat = "general"
$ws.Range("I43").Value = "13.1 mm"
$ws.Range("E44").Value = "2026-02-13 19:50:42"
$ws.Range("I44").Value = "6.9 mm"
$ws.Range("E45").Value = "2026-02-13 19:50:45"
$ws.Range("H45").NumberFormat = "@"
$ws.Range("H45").Value = "64%"
$ws.Range("H45").NumberFormat = "general"
$ws.Range("J45").Value = "993.4 hPa"
$ws.Range("O45").Value = "5.7 °C"
$ws.Range("E46").Value = "2026-02-13 19:50:47"
$ws.Range("J46").Value = "995.6 hPa"
$ws.Range("O46").Value = "8.9 °C"
